$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 159
$ws.Range("F4").Value = 1093
$ws.Range("F6").Value = 340
$ws.Range("F8").Value = 12592
$ws.Range("F9").Value = 2217
$ws.Range("F10").Value = 909
$ws.Range("F12").Value = 23810
$ws.Range("F14").Value = 1254
$ws.Range("F15").Value = 244
$ws.Range("F16").Value = 284
$ws.Range("F18").Value = 685
$ws.Range("F21").Value = 791
$ws.Range("F22").Value = 4518
$ws.Range("F23").Value = 1165
$ws.Range("F24").Value = 891
$ws.Range("F29").Value = 1111
$ws.Range("F31").Value = 125
$ws.Range("F32").Value = 287
$ws.Range("F36").Value = 29
$ws.Range("F37").Value = 4517
$ws.Range("F39").Value = 4628
$ws.Range("F40").Value = 5589
$ws.Range("F45").Value = 374
$ws.Range("F48").Value = 4125
$ws.Range("F49").Value = 153

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 106
$ws.Range("F12").Value = 1063
$ws.Range("F24").Value = 8

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 767
$ws.Range("F4").Value = 111
$ws.Range("F5").Value = 22

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 159
$ws.Range("F5").Value = 1093
$ws.Range("F7").Value = 340
$ws.Range("F9").Value = 12592
$ws.Range("F10").Value = 2217
$ws.Range("F11").Value = 909
$ws.Range("F13").Value = 1254
$ws.Range("F14").Value = 244
$ws.Range("F15").Value = 284
$ws.Range("F17").Value = 685
$ws.Range("F20").Value = 791
$ws.Range("F21").Value = 4518
$ws.Range("F22").Value = 4518
$ws.Range("F23").Value = 1165
$ws.Range("F24").Value = 22
$ws.Range("F25").Value = 106
$ws.Range("F31").Value = 1111
$ws.Range("F33").Value = 125
$ws.Range("F35").Value = 287
$ws.Range("F39").Value = 4628
$ws.Range("F46").Value = 4125
$ws.Range("F50").Value = 8
